$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated experiment results (2021/04/21 15:43) ---
$ws.Range("B24").Value = 0.83160000000000001
$ws.Range("B26").Value = 0.89739999999999998
$ws.Range("C26").Value = 37
$ws.Range("B27").Value = 0.88009999999999999
$ws.Range("C27").Value = 38

# B32 (=AVERAGE(B2:B31)) recalculates automatically from the edits above.

# The "平均" label (A32) and its average value (B32) keep the same visual
# formatting (A32 right-aligned, B32 default) but the underlying style
# records end up swapped relative to each other after the resave - swap
# the two cells' formats to mirror that.
$ws.Range("A32").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4122) | Out-Null
$ws.Range("B32").Copy() | Out-Null
$ws.Range("E32").PasteSpecial(-4122) | Out-Null

$ws.Range("A32").ClearFormats()
$ws.Range("B32").ClearFormats()

$ws.Range("E32").Copy() | Out-Null
$ws.Range("A32").PasteSpecial(-4122) | Out-Null
$ws.Range("D32").Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4122) | Out-Null

$ws.Range("D32:E32").Clear()
